$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 275, shifting existing rows 275-368 down to 276-369
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new record
$ws.Range("A275").Value = 9
$ws.Range("B275").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C275").Value = "Metropolitana"
$ws.Range("D275").Value = 44988
$ws.Range("E275").Value = 13
$ws.Range("F275").Value = 100112030
$ws.Range("G275").Value = "Poroto granado"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 70
$ws.Range("K275").Value = 25000
$ws.Range("L275").Value = 27000
$ws.Range("M275").Value = 26000
$ws.Range("N275").Value = "$/saco 25 kilos"
$ws.Range("O275").Value = "Región Metropolitana"
$ws.Range("P275").Value = 1040
$ws.Range("Q275").Value = 25
$ws.Range("R275").Value = "Hortaliza"
